$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A for the tab-name labels.
$ws.Columns("A:A").Insert()

# New "TabName" / "CasesTab" column content.
$ws.Range("A1").Value() = "TabName"
$ws.Range("A2").Value() = "CasesTab"

# Updated Cypher query text for the Cases tab (now in column B, row 2).
$casesQuery = @'
MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = "NOT_REPORTED"
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS `Case ID`,
    COALESCE(ct.clinical_trial_designation, '') AS `Trial Code`,
    COALESCE(a.arm_id, '') AS `Arm`,
    COALESCE(a.arm_drug, '') AS `Arm Treatment`,
    COALESCE(c.disease, '') AS `Diagnosis`,
    COALESCE(c.gender, '') AS `Gender`,
    COALESCE(c.race, '') AS `Race`,
    COALESCE(c.ethnicity, '') AS `Ethnicity`
'@

# Updated Cypher query text for the Stat tab (now in column C, row 2).
$statQuery = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = "NOT_REPORTED"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials
'@

$ws.Range("B2").Value() = $casesQuery
$ws.Range("C2").Value() = $statQuery

# C2 (the stat query) now needs the same wrap-text style already used by B2.
$ws.Range("C2").WrapText = $true

# Row 2 grew taller to accommodate the longer wrapped text.
$ws.Rows(2).RowHeight = 174

# Column widths: new column A is narrow, others keep/shift their prior widths.
$ws.Columns("A:A").ColumnWidth = 8

# Update dimension/selection bookkeeping to match the new layout.
$ws.Range("B2").Select()
